# Applies the OSA1 homework content edits described in the commit.
# The surrounding diff is dominated by Word's proofing engine re-splitting
# runs (proofErr/gramStart/gramEnd markers) and repagination
# (lastRenderedPageBreak) artifacts that carry no visible text change; the
# substantive wording changes are the two replacements below.

$d = $word.ActiveDocument

# 1) Clarify that it's *instructions* (not whole programs) that run at a
#    time on a single-core stored-program machine.
$d.Content.Find.Execute(
    "if the OS is running, then a program is not. Similarly, if a program is running, then neither the OS nor any other program is running.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "if an instruction in the OS is running, then an instruction from a program is not. Similarly, if an instruction from a program is running, then neither instructions from the OS nor any other program are running.",
    2
) | Out-Null

# 2) Rework the "summarize this tradeoff" prompt to ask how the
#    convenience/efficiency tradeoff differs across device classes.
$d.Content.Find.Execute(
    " with which application programs are executed.  To increase user convenience, we must decrease the efficiency of user programs or vice versa. Summarize this tradeoff and why it exists in a few sentences of your own words.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " with which application programs are executed.  In a few sentences of your own words, summarize this tradeoff and how it might be made differently in an OS for a smart watch, a personal computer, and a high performance super computer.",
    2
) | Out-Null
